{"js": "const replacements = [\n  [\"2024-01-30 Tuesday\", \"2024-01-31 Wednesday\"],\n  [\"86\u00d732=\", \"45\u00d754=\"],\n  [\"93\u00d776=\", \"64\u00d718=\"],\n  [\"21\u00d747=\", \"57\u00d798=\"],\n  [\"11\u00d775=\", \"33\u00d744=\"],\n  [\"63\u00d729=\", \"59\u00d757=\"],\n  [\"51\u00d780=\", \"40\u00d724=\"],\n  [\"22\u00d721=\", \"50\u00d737=\"],\n  [\"45\u00d739=\", \"88\u00d753=\"],\n  [\"88\u00d790=\", \"28\u00d792=\"],\n  [\"34\u00d772=\", \"60\u00d782=\"],\n  [\"25\u00d725=\", \"12\u00d771=\"],\n  [\"58\u00d728=\", \"95\u00d754=\"],\n  [\"71\u00d750=\", \"25\u00d769=\"],\n  [\"57\u00d782=\", \"75\u00d774=\"],\n  [\"68\u00d799=\", \"19\u00d718=\"],\n  [\"76\u00d752=\", \"73\u00d727=\"],\n  [\"64\u00d763=\", \"64\u00d719=\"],\n  [\"11\u00d742=\", \"80\u00d749=\"],\n  [\"98\u00d774=\", \"25\u00d777=\"],\n  [\"77\u00d795=\", \"63\u00d761=\"],\n  [\"22\u00d736=\", \"74\u00d784=\"],\n  [\"71\u00d729=\", \"68\u00d769=\"],\n  [\"44\u00d723=\", \"24\u00d785=\"],\n  [\"46\u00d713=\", \"74\u00d762=\"],\n  [\"80\u00d783=\", \"11\u00d716=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n", "ps1": "$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"2024-01-30 Tuesday\"\n$find.Replacement.Text = \"2024-01-31 Wednesday\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"86\u00d732=\"\n$find.Replacement.Text = \"45\u00d754=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"93\u00d776=\"\n$find.Replacement.Text = \"64\u00d718=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"21\u00d747=\"\n$find.Replacement.Text = \"57\u00d798=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"11\u00d775=\"\n$find.Replacement.Text = \"33\u00d744=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"63\u00d729=\"\n$find.Replacement.Text = \"59\u00d757=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"51\u00d780=\"\n$find.Replacement.Text = \"40\u00d724=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"22\u00d721=\"\n$find.Replacement.Text = \"50\u00d737=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"45\u00d739=\"\n$find.Replacement.Text = \"88\u00d753=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"88\u00d790=\"\n$find.Replacement.Text = \"28\u00d792=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"34\u00d772=\"\n$find.Replacement.Text = \"60\u00d782=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"25\u00d725=\"\n$find.Replacement.Text = \"12\u00d771=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"58\u00d728=\"\n$find.Replacement.Text = \"95\u00d754=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"71\u00d750=\"\n$find.Replacement.Text = \"25\u00d769=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"57\u00d782=\"\n$find.Replacement.Text = \"75\u00d774=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"68\u00d799=\"\n$find.Replacement.Text = \"19\u00d718=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"76\u00d752=\"\n$find.Replacement.Text = \"73\u00d727=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"64\u00d763=\"\n$find.Replacement.Text = \"64\u00d719=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"11\u00d742=\"\n$find.Replacement.Text = \"80\u00d749=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"98\u00d774=\"\n$find.Replacement.Text = \"25\u00d777=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"77\u00d795=\"\n$find.Replacement.Text = \"63\u00d761=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"22\u00d736=\"\n$find.Replacement.Text = \"74\u00d784=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"71\u00d729=\"\n$find.Replacement.Text = \"68\u00d769=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"44\u00d723=\"\n$find.Replacement.Text = \"24\u00d785=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"46\u00d713=\"\n$find.Replacement.Text = \"74\u00d762=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"80\u00d783=\"\n$find.Replacement.Text = \"11\u00d716=\"\n$find.Execute($find.Text, $false, $true, $false, $false, $false, $true, 1, $true, $find.Replacement.Text, 2)\n\n"}
